# Auto-generated: apply scheduled-runner market data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6994179
$ws.Range("I112").Value = 716.6667
$ws.Range("J112").Value = 8265717.5
$ws.Range("K112").Value = 2150.0001
$ws.Range("L112").Value = 24797152.5
$ws.Range("M112").Value = -1042.0001
$ws.Range("N112").Value = -24799368.5

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = None
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null

$ws.Range("H116").Value = 6294330
$ws.Range("I116").Value = 9889828
$ws.Range("J116").Value = 2208.25
$ws.Range("K116").Value = 9889828
$ws.Range("L116").Value = 2208.25
$ws.Range("M116").Value = -9886386
$ws.Range("N116").Value = -9092.25

$ws.Range("H132").Value = 195332.27
$ws.Range("I132").Value = 217793.34
$ws.Range("J132").Value = 15643.714
$ws.Range("K132").Value = 653380.02
$ws.Range("L132").Value = 46931.142
$ws.Range("M132").Value = -650850.02
$ws.Range("N132").Value = -51991.142

$ws.Range("H137").Value = 20000920
$ws.Range("I137").Value = 26316418
$ws.Range("J137").Value = 1836.4166
$ws.Range("K137").Value = 78949254
$ws.Range("L137").Value = 5509.2498
$ws.Range("M137").Value = -78946704
$ws.Range("N137").Value = -10609.2498

$ws.Range("H138").Value = 3907226.8
$ws.Range("I138").Value = 1198325.1
$ws.Range("J138").Value = 5651943
$ws.Range("K138").Value = 3594975.3
$ws.Range("L138").Value = 16955829
$ws.Range("M138").Value = -3589835.3
$ws.Range("N138").Value = -16966109

$ws.Range("H141").Value = 2000.4066
$ws.Range("J141").Value = 6353.2666
$ws.Range("L141").Value = 19059.7998
$ws.Range("N141").Value = -29419.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1093.21
$ws.Range("I32").Value = 1029.4517
$ws.Range("K32").Value = 1029.4517
$ws.Range("M32").Value = -742.4517000000001

$ws.Range("H61").Value = 2248.8865
$ws.Range("I61").Value = 1686.8918
$ws.Range("J61").Value = 5219.4287
$ws.Range("K61").Value = 1686.8918
$ws.Range("L61").Value = 5219.4287
$ws.Range("M61").Value = -1474.8918
$ws.Range("N61").Value = -5643.4287

$ws.Range("H74").Value = 5377.5835
$ws.Range("I74").Value = 1810.3334
$ws.Range("J74").Value = 16079.333
$ws.Range("K74").Value = 1810.3334
$ws.Range("L74").Value = 16079.333
$ws.Range("M74").Value = -936.3334
$ws.Range("N74").Value = -17827.333

$ws.Range("H77").Value = 5377.5835
$ws.Range("I77").Value = 1810.3334
$ws.Range("J77").Value = 16079.333
$ws.Range("K77").Value = 9051.666999999999
$ws.Range("L77").Value = 80396.66500000001
$ws.Range("M77").Value = -4683.666999999999
$ws.Range("N77").Value = -89132.66500000001

$ws.Range("H132").Value = 1614
$ws.Range("I132").Value = 1158.0426
$ws.Range("J132").Value = 4292.75
$ws.Range("K132").Value = 3474.1278
$ws.Range("L132").Value = 12878.25
$ws.Range("M132").Value = -944.1278000000002
$ws.Range("N132").Value = -17938.25

$ws.Range("H133").Value = 53666.332
$ws.Range("J133").Value = 53666.332
$ws.Range("L133").Value = 53666.332
$ws.Range("N133").Value = -58726.332

$ws.Range("H136").Value = 2248.8865
$ws.Range("I136").Value = 1686.8918
$ws.Range("J136").Value = 5219.4287
$ws.Range("K136").Value = 5060.6754
$ws.Range("L136").Value = 15658.2861
$ws.Range("M136").Value = -2510.6754
$ws.Range("N136").Value = -20758.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1133.2115
$ws.Range("I58").Value = 760.73914
$ws.Range("J58").Value = 3988.8333
$ws.Range("K58").Value = 760.73914
$ws.Range("L58").Value = 3988.8333
$ws.Range("M58").Value = -557.73914
$ws.Range("N58").Value = -4394.8333

$ws.Range("H132").Value = 1309.5682
$ws.Range("I132").Value = 1063.975
$ws.Range("J132").Value = 3765.5
$ws.Range("K132").Value = 3191.925
$ws.Range("L132").Value = 11296.5
$ws.Range("M132").Value = -661.9249999999997
$ws.Range("N132").Value = -16356.5

$ws.Range("H134").Value = 1481.131
$ws.Range("I134").Value = 964.24677
$ws.Range("J134").Value = 7166.857
$ws.Range("K134").Value = 2892.74031
$ws.Range("L134").Value = 21500.571
$ws.Range("M134").Value = -357.7403100000001
$ws.Range("N134").Value = -26570.571

$ws.Range("H136").Value = 1133.2115
$ws.Range("I136").Value = 760.73914
$ws.Range("J136").Value = 3988.8333
$ws.Range("K136").Value = 2282.21742
$ws.Range("L136").Value = 11966.4999
$ws.Range("M136").Value = 267.7825800000001
$ws.Range("N136").Value = -17066.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1198.3334
$ws.Range("J80").Value = 1198.3334
$ws.Range("L80").Value = 3595.0002
$ws.Range("N80").Value = -5467.0002

$ws.Range("H83").Value = 1198.3334
$ws.Range("J83").Value = 1198.3334
$ws.Range("L83").Value = 10785.0006
$ws.Range("N83").Value = -20145.0006

$ws.Range("H92").Value = 789.44446
$ws.Range("J92").Value = 780.6
$ws.Range("L92").Value = 2341.8
$ws.Range("N92").Value = -4837.8

$ws.Range("H129").Value = 1539.1538
$ws.Range("I129").Value = 1440
$ws.Range("J129").Value = 1591.6471
$ws.Range("K129").Value = 4320
$ws.Range("L129").Value = 4774.9413
$ws.Range("M129").Value = 680
$ws.Range("N129").Value = -14774.9413

$ws.Range("H132").Value = 2683.8333
$ws.Range("I132").Value = 2597.5
$ws.Range("K132").Value = 23377.5
$ws.Range("M132").Value = -20847.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 619020.4399999999
$ws.Range("I122").Value = 742224
$ws.Range("J122").Value = 3002.6667
$ws.Range("K122").Value = 2226672
$ws.Range("L122").Value = 9008.000100000001
$ws.Range("M122").Value = -2224222
$ws.Range("N122").Value = -13908.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6067.727
$ws.Range("I22").Value = 966.6667
$ws.Range("J22").Value = 7980.625
$ws.Range("K22").Value = 966.6667
$ws.Range("L22").Value = 7980.625
$ws.Range("M22").Value = -671.6667
$ws.Range("N22").Value = -8570.625

$ws.Range("H27").Value = 6067.727
$ws.Range("I27").Value = 966.6667
$ws.Range("J27").Value = 7980.625
$ws.Range("K27").Value = 966.6667
$ws.Range("L27").Value = 7980.625
$ws.Range("M27").Value = -859.6667
$ws.Range("N27").Value = -8194.625

$ws.Range("H122").Value = 3348
$ws.Range("I122").Value = 2915.4285
$ws.Range("J122").Value = 3564.2856
$ws.Range("K122").Value = 8746.2855
$ws.Range("L122").Value = 10692.8568
$ws.Range("M122").Value = -6296.2855
$ws.Range("N122").Value = -15592.8568

$ws.Range("H132").Value = 3557.3188
$ws.Range("I132").Value = 2855.9614
$ws.Range("J132").Value = 5702.647
$ws.Range("K132").Value = 8567.8842
$ws.Range("L132").Value = 17107.941
$ws.Range("M132").Value = -6037.8842
$ws.Range("N132").Value = -22167.941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6580640.5
$ws.Range("I132").Value = 8773318
$ws.Range("J132").Value = 2609
$ws.Range("K132").Value = 26319954
$ws.Range("L132").Value = 7827
$ws.Range("M132").Value = -26317424
$ws.Range("N132").Value = -12887

$ws.Range("H136").Value = 19997.227
$ws.Range("I136").Value = 24977.39
$ws.Range("J136").Value = 2981.6667
$ws.Range("K136").Value = 74932.17
$ws.Range("L136").Value = 8945.000100000001
$ws.Range("M136").Value = -72382.17
$ws.Range("N136").Value = -14045.0001
